# Lab5: extensor de sinal, 4 novos estados na UC e registrador de instrução
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (header row) ---
$ws.Range("A3").Value = "Nome"
$ws.Range("B3").Value = "formato"
$ws.Range("C3").Value = "Codigo"
$ws.Range("D3").Value = "vetor"

# --- Row 4 (Jump) ---
$ws.Range("A4").Value = "Jump"
$ws.Range("B4").Value = "J"
$ws.Range("C4").Value = "func:00,op:111"
$ws.Range("D4").Value = "[reserv17:12|end11:5|func4:3|op2:0]"

# --- Row 5 (Nop) ---
$ws.Range("A5").Value = "Nop"
$ws.Range("B5").Value = "N"
$ws.Range("C5").Value = "func:00,op:000"
$ws.Range("D5").Value = "[zeros17:5|func4:3|op2:0]"

# --- Row 6 (Ld - acc) ---
$ws.Range("A6").Value = "Ld"
$ws.Range("B6").Value = "acc <- const"
$ws.Range("C6").Value = "func:00,op:001"
$ws.Range("D6").Value = "[const17:11|reservado10:6|acc5|func4:3|op2:0]"

# --- Row 7 (Ld - Rn) ---
$ws.Range("A7").Value = "Ld"
$ws.Range("B7").Value = "Rn <- const"
$ws.Range("C7").Value = "func:01,op:001"
$ws.Range("D7").Value = "[const17:11|reservado10:8|rn7:5|func4:3|op2:0]"

# --- Row 8 (Mov A<-Rn) ---
$ws.Range("A8").Value = "Mov"
$ws.Range("B8").Value = "A<-Rn"
$ws.Range("C8").Value = "func:00,op:010"
$ws.Range("D8").Value = "[const17:11?|reservado10:9|rn8:6|acc5|func4:3|op2:0]"

# --- Row 9 (Mov Rn<-A) ---
$ws.Range("A9").Value = "Mov"
$ws.Range("B9").Value = "Rn<-A"
$ws.Range("C9").Value = "func:01,op:010"
$ws.Range("D9").Value = "[const17:11|reservado10:9|rn8:6|acc5|func4:3|op2:0]"
$ws.Range("D9").Font.Underline = $true

# --- Row 10 (op(tipo) - soma) ---
$ws.Range("A10").Value = "op(tipo)"
$ws.Range("B10").Value = "Rn <-A"
$ws.Range("C10").Value = "func:00 soma op:011"
$ws.Range("D10").Value = "[const17:11|tipo10:9|rn8:6|acc5|func4:3|op2:0]"

# --- Row 11 (op(tipo) - Subtração) ---
$ws.Range("A11").Value = "op(tipo)"
$ws.Range("B11").Value = "Rn <-A"
$ws.Range("C11").Value = "func:01 Subtração op:011"
$ws.Range("D11").Value = "[const17:11|tipo10:9|rn8:6|acc5|func4:3|op2:0]"

# --- Row 12 (op(tipo) - Nand) ---
$ws.Range("A12").Value = "op(tipo)"
$ws.Range("B12").Value = "Rn <-A"
$ws.Range("C12").Value = "func:00 Nand op:011"
$ws.Range("D12").Value = "[const17:11|tipo10:9|rn8:6|acc5|func4:3|op2:0]"

# --- Row 13 (op(tipo) - Xor) ---
$ws.Range("A13").Value = "op(tipo)"
$ws.Range("B13").Value = "Rn <-A"
$ws.Range("C13").Value = "func:00 Xor op:011"
$ws.Range("D13").Value = "[const17:11|tipo10:9|rn8:6|acc5|func4:3|op2:0]"

# --- Column D width widened to fit the longer vetor strings ---
$ws.Columns("D").ColumnWidth = 49

# --- View: scroll down a bit and select D9 ---
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("D9").Select() | Out-Null

# --- Page setup (A4 portrait) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
